$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function RGBColor($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }
$green  = RGBColor 0 176 80    # FF00B050 - "done"
$orange = RGBColor 255 192 0   # FFFFC000 - "in progress"

# Row 5 ("find out how much each measurement occurs...") - mark Status green (done)
# and add a comment describing where things stand.
$ws.Range("D5").Interior.Color = $green
$ws.Range("E5").Value = "Visualization is done" + [char]10 + "The only thing left is to decide which measurements to discard "
$ws.Range("E5").WrapText = $true

# Row 6 ("If >= 2 data points: fit (linear/quadr./...) for measurements...")
# mark Status orange (started) with matching orange text, plus a comment.
$ws.Range("D6").Font.Color = $orange
$ws.Range("D6").Interior.Color = $orange
$ws.Range("E6").Value = "Started, has to be improved"

# Row 8 ("If all nan: either from other measurements, or with global mean.")
# mark Status orange (started) plus a comment.
$ws.Range("D8").Interior.Color = $orange
$ws.Range("E8").Value = "global median for now"

# Reflect the author's final cursor position/selection when the file was saved
$ws.Range("E10").Select() | Out-Null
